# Add a new "2022-Q3" sheet (duplicated from "2022-Q2" to inherit the header
# styling/layout) positioned right after "总计" and before "2022-Q2", then
# overwrite its data with the 2022-Q3 fund figures. Finally, insert a
# matching summary row on the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q3" worksheet by copying "2022-Q2" (Excel places
#    the copy immediately before the source sheet), then rename it and trim
#    it down to the 6 data rows required for 2022-Q3 (source has 7).
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"
$q3.Rows.Item(8).Delete()

# ---------------------------------------------------------------------
# 2) Overwrite the data rows (2-7) with the 2022-Q3 fund holdings.
#    Columns: B=code, C=name, D=fund size, E=stock position, F=position pct,
#    G=held value (100M), H=position rank. B and D-G are text in the source
#    file (e.g. fund codes keep leading zeros); force text storage so the
#    numeric-looking strings don't get reinterpreted as numbers.
# ---------------------------------------------------------------------
$q3.Range("B2:B7").NumberFormat = "@"
$q3.Range("D2:G7").NumberFormat = "@"

$q3.Range("B2").Value = "005775"
$q3.Range("C2").Value = "中加转型动力灵活配置混合A"
$q3.Range("D2").Value = "5.34"
$q3.Range("E2").Value = "50.55"
$q3.Range("F2").Value = "2.37"
$q3.Range("G2").Value = "0.1266"
$q3.Range("H2").Value = 10

$q3.Range("B3").Value = "009242"
$q3.Range("C3").Value = "中加核心智造混合A"
$q3.Range("D3").Value = "1.92"
$q3.Range("E3").Value = "61.20"
$q3.Range("F3").Value = "2.83"
$q3.Range("G3").Value = "0.0543"
$q3.Range("H3").Value = 6

$q3.Range("B4").Value = "012072"
$q3.Range("C4").Value = "中加喜利回报一年持有期混合C"
$q3.Range("D4").Value = "2.21"
$q3.Range("E4").Value = "38.64"
$q3.Range("F4").Value = "2.04"
$q3.Range("G4").Value = "0.0451"
$q3.Range("H4").Value = 7

$q3.Range("B5").Value = "012071"
$q3.Range("C5").Value = "中加喜利回报一年持有期混合A"
$q3.Range("D5").Value = "1.98"
$q3.Range("E5").Value = "38.64"
$q3.Range("F5").Value = "2.04"
$q3.Range("G5").Value = "0.0404"
$q3.Range("H5").Value = 7

$q3.Range("B6").Value = "005776"
$q3.Range("C6").Value = "中加转型动力灵活配置混合C"
$q3.Range("D6").Value = "1.27"
$q3.Range("E6").Value = "50.55"
$q3.Range("F6").Value = "2.37"
$q3.Range("G6").Value = "0.0301"
$q3.Range("H6").Value = 10

$q3.Range("B7").Value = "009243"
$q3.Range("C7").Value = "中加核心智造混合C"
$q3.Range("D7").Value = "0.09"
$q3.Range("E7").Value = "61.20"
$q3.Range("F7").Value = "2.83"
$q3.Range("G7").Value = "0.0025"
$q3.Range("H7").Value = 6

# ---------------------------------------------------------------------
# 3) Update the "总计" (summary) sheet: the quarter/count/value columns
#    (B:D) shift down by one data row to make room for 2022-Q3 at the top,
#    while column A keeps its existing running index (0,1,2,...) per row
#    position and simply gains one more entry (6) for the newly appended
#    row 8. Values are written bottom-up so a row's old data is read
#    before it gets overwritten.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

# New row 8 (was not present before) - clone A7's style onto A8 first.
$total.Range("A7").Copy()
$total.Range("A8").PasteSpecial(-4122)
$total.Range("A8").Value = 6
$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 8
$total.Range("D8").Value = 0.24

$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 9
$total.Range("D7").Value = 0.88

$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 4
$total.Range("D6").Value = 0.33

$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 13
$total.Range("D5").Value = 0.92

$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 6
$total.Range("D4").Value = 0.5600000000000001

$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 7
$total.Range("D3").Value = 0.83

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.3

$q3.Range("A1").Select()
